$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value = 1510.3125
$ws.Range("J17").Value = 1510.3125
$ws.Range("L17").Value = 4530.9375
$ws.Range("N17").Value = -4866.9375

$ws.Range("H28").Value = 1122.85
$ws.Range("I28").Value = 190.45454
$ws.Range("J28").Value = 2262.4443
$ws.Range("K28").Value = 190.45454
$ws.Range("L28").Value = 2262.4443
$ws.Range("M28").Value = 294.54546
$ws.Range("N28").Value = -3232.4443

$ws.Range("H98").Value = 1279.9231
$ws.Range("I98").Value = 1279.9231
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1279.9231
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 218.0769
$ws.Range("N98").ClearContents()

$ws.Range("H100").Value = 1128.8889
$ws.Range("I100").Value = 1131.75
$ws.Range("J100").Value = 1106
$ws.Range("K100").Value = 1131.75
$ws.Range("L100").Value = 1106
$ws.Range("M100").Value = -590.75
$ws.Range("N100").Value = -2188

$ws.Range("H122").Value = 1279.9231
$ws.Range("I122").Value = 1279.9231
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3839.7693
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1389.7693
$ws.Range("N122").ClearContents()

$ws.Range("H137").Value = 1295.3549
$ws.Range("I137").Value = 908.36365
$ws.Range("J137").Value = 2241.3333
$ws.Range("K137").Value = 2725.09095
$ws.Range("L137").Value = 6723.999899999999
$ws.Range("M137").Value = -175.0909499999998
$ws.Range("N137").Value = -11823.9999

$ws = $wb.Worksheets.Item(2)
$ws.Range("H102").Value = 1999.5
$ws.Range("I102").Value = 1999
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1999
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -377
$ws.Range("N102").Value = -5244

$ws.Range("H122").Value = 1529
$ws.Range("I122").Value = 1529
$ws.Range("K122").Value = 4587
$ws.Range("M122").Value = -2137

$ws.Range("H132").Value = 4964.811
$ws.Range("I132").Value = 5055
$ws.Range("J132").Value = 4498.8335
$ws.Range("K132").Value = 15165
$ws.Range("L132").Value = 13496.5005
$ws.Range("M132").Value = -12635
$ws.Range("N132").Value = -18556.5005

$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 4025.6875
$ws.Range("I86").Value = 2823.5557
$ws.Range("J86").Value = 5571.2856
$ws.Range("K86").Value = 2823.5557
$ws.Range("L86").Value = 5571.2856
$ws.Range("M86").Value = -1700.5557
$ws.Range("N86").Value = -7817.2856

$ws.Range("H89").Value = 4025.6875
$ws.Range("I89").Value = 2823.5557
$ws.Range("J89").Value = 5571.2856
$ws.Range("K89").Value = 14117.7785
$ws.Range("L89").Value = 27856.428
$ws.Range("M89").Value = -8501.7785
$ws.Range("N89").Value = -39088.428

$ws.Range("H94").Value = 899.8
$ws.Range("I94").Value = 974.75
$ws.Range("J94").Value = 600
$ws.Range("K94").Value = 974.75
$ws.Range("L94").Value = 600
$ws.Range("M94").Value = -523.75
$ws.Range("N94").Value = -1502

$ws.Range("H105").Value = 4436
$ws.Range("I105").Value = 4542.5713
$ws.Range("J105").Value = 4249.5
$ws.Range("K105").Value = 4542.5713
$ws.Range("L105").Value = 4249.5
$ws.Range("M105").Value = -2795.5713
$ws.Range("N105").Value = -7743.5

$ws.Range("H107").Value = 875.1667
$ws.Range("I107").Value = 650.2
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 650.2
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 1269.8
$ws.Range("N107").Value = -5840

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 2566086.8
$ws.Range("I31").Value = 1489.5682
$ws.Range("J31").Value = 7939528.5
$ws.Range("K31").Value = 1489.5682
$ws.Range("L31").Value = 7939528.5
$ws.Range("M31").Value = -1194.5682
$ws.Range("N31").Value = -7940118.5

$ws.Range("H34").Value = 2566086.8
$ws.Range("I34").Value = 1489.5682
$ws.Range("J34").Value = 7939528.5
$ws.Range("K34").Value = 1489.5682
$ws.Range("L34").Value = 7939528.5
$ws.Range("M34").Value = -1287.5682
$ws.Range("N34").Value = -7939932.5

$ws.Range("H86").Value = 2386
$ws.Range("I86").Value = 1182.7142
$ws.Range("J86").Value = 3228.3
$ws.Range("K86").Value = 1182.7142
$ws.Range("L86").Value = 3228.3
$ws.Range("M86").Value = -59.71419999999989
$ws.Range("N86").Value = -5474.3

$ws.Range("H89").Value = 2386
$ws.Range("I89").Value = 1182.7142
$ws.Range("J89").Value = 3228.3
$ws.Range("K89").Value = 5913.571
$ws.Range("L89").Value = 16141.5
$ws.Range("M89").Value = -297.5709999999999
$ws.Range("N89").Value = -27373.5

$ws.Range("H99").Value = 2243.0557
$ws.Range("I99").Value = 2137
$ws.Range("J99").Value = 2349.111
$ws.Range("K99").Value = 2137
$ws.Range("L99").Value = 2349.111
$ws.Range("M99").Value = -639
$ws.Range("N99").Value = -5345.111

$ws.Range("H126").Value = 2243.0557
$ws.Range("I126").Value = 2137
$ws.Range("J126").Value = 2349.111
$ws.Range("K126").Value = 6411
$ws.Range("L126").Value = 7047.333
$ws.Range("M126").Value = -3941
$ws.Range("N126").Value = -11987.333

$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 2367167.2
$ws.Range("I131").Value = 8283.214
$ws.Range("J131").Value = 3367905.8
$ws.Range("K131").Value = 24849.642
$ws.Range("L131").Value = 10103717.4
$ws.Range("M131").Value = -19809.642
$ws.Range("N131").Value = -10113797.4

$ws = $wb.Worksheets.Item(6)
$ws.Range("H45").Value = 36000
$ws.Range("J45").Value = 36000
$ws.Range("L45").Value = 36000
$ws.Range("N45").Value = -37118

$ws.Range("H51").Value = 35500
$ws.Range("J51").Value = 35500
$ws.Range("L51").Value = 35500
$ws.Range("N51").Value = -36518

$ws.Range("H80").Value = 5271.357
$ws.Range("I80").Value = 4757
$ws.Range("J80").Value = 5785.7144
$ws.Range("K80").Value = 4757
$ws.Range("L80").Value = 5785.7144
$ws.Range("M80").Value = -3759
$ws.Range("N80").Value = -7781.7144

$ws.Range("H83").Value = 5271.357
$ws.Range("I83").Value = 4757
$ws.Range("J83").Value = 5785.7144
$ws.Range("K83").Value = 23785
$ws.Range("L83").Value = 28928.572
$ws.Range("M83").Value = -18793
$ws.Range("N83").Value = -38912.572

$ws.Range("H102").Value = 1800
$ws.Range("I102").Value = 1750
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1750
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -128
$ws.Range("N102").Value = -5244

$ws.Range("H126").Value = 4230.75
$ws.Range("I126").Value = 5303
$ws.Range("J126").Value = 1014
$ws.Range("K126").Value = 15909
$ws.Range("L126").Value = 3042
$ws.Range("M126").Value = -13439
$ws.Range("N126").Value = -7982

$ws = $wb.Worksheets.Item(7)
$ws.Range("H68").Value = 1746.25
$ws.Range("I68").Value = 1371.1111
$ws.Range("J68").Value = 2228.5715
$ws.Range("K68").Value = 1371.1111
$ws.Range("L68").Value = 2228.5715
$ws.Range("M68").Value = -622.1111000000001
$ws.Range("N68").Value = -3726.5715

$ws.Range("H71").Value = 1746.25
$ws.Range("I71").Value = 1371.1111
$ws.Range("J71").Value = 2228.5715
$ws.Range("K71").Value = 6855.5555
$ws.Range("L71").Value = 11142.8575
$ws.Range("M71").Value = -3111.5555
$ws.Range("N71").Value = -18630.8575

$ws.Range("H93").Value = 1983.9166
$ws.Range("I93").Value = 1980.6
$ws.Range("J93").Value = 1986.2858
$ws.Range("K93").Value = 1980.6
$ws.Range("L93").Value = 1986.2858
$ws.Range("M93").Value = -732.5999999999999
$ws.Range("N93").Value = -4482.2858

$ws.Range("H100").Value = 1710.2
$ws.Range("I100").Value = 1700.3572
$ws.Range("J100").Value = 1733.1666
$ws.Range("K100").Value = 1700.3572
$ws.Range("L100").Value = 1733.1666
$ws.Range("M100").Value = -1159.3572
$ws.Range("N100").Value = -2815.1666

$ws.Range("H133").Value = 19665.2
$ws.Range("J133").Value = 19665.2
$ws.Range("L133").Value = 19665.2
$ws.Range("N133").Value = -24725.2

$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H96").Value = 62500624
$ws.Range("I96").Value = 250000000
$ws.Range("J96").Value = 833.3333
$ws.Range("K96").Value = 250000000
$ws.Range("L96").Value = 833.3333
$ws.Range("M96").Value = -249998627
$ws.Range("N96").Value = -3579.3333

$ws.Range("H122").Value = 46344.523
$ws.Range("I122").Value = 65310
$ws.Range("J122").Value = 2994.8572
$ws.Range("K122").Value = 195930
$ws.Range("L122").Value = 8984.571599999999
$ws.Range("M122").Value = -193480
$ws.Range("N122").Value = -13884.5716

$ws.Range("H127").Value = 20429
$ws.Range("J127").Value = 20429
$ws.Range("L127").Value = 20429
$ws.Range("N127").Value = -30349
